$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra rows (3 through 13), keep only header + one data row.
$ws.Rows("3:13").Delete()

# Update the remaining data row values to the new realistic demo data.
# (set B2 first so the shared-string table order matches: ofs-pop before population)
$ws.Range("B2").Value = "ofs-pop"
$ws.Range("A2").Value = "population"

# Resize the Excel table to the new data extent.
$table = $ws.ListObjects.Item("Tableau1")
$table.Resize($ws.Range("A1:B2"))

# Adjust column A width to fit new (longer) content ("population").
$ws.Columns("A").ColumnWidth = 9

# Update the active selection in the sheet view to match saved state.
$ws.Range("B6").Select()
